# LORILLA, LOIDA.xlsx -- "Add Leave Card 10/32023 3:18 PM"
#
# Summary of the edit being replayed:
#  1. On "2018 LEAVE CREDITS" (sheet1 / Table13), the monthly earn
#     calculator (col A/C/G) gains 4 more months (Aug-Nov 2023) and the
#     leave-usage annotations that used to live in rows 82-84
#     (SL(2-0-0), SL(4-0-0), SP(1-0-0)) are cleared out of that sheet.
#  2. Those leave-usage annotations (SL(2-0-0) & SL(4-0-0), with their
#     day counts and remark dates) are recreated on "2017 LEAVE BALANCE"
#     (sheet2 / Table1) rows 18-19, and two brand new VL(8-0-0) entries
#     are appended on rows 20-21 for September & October 2023.
#  3. Two now-unused trailing blank rows (125 & 126) are deleted from the
#     2018 LEAVE CREDITS table, shrinking it from A8:K127 to A8:K125.
#  4. The active sheet/tab switches from "2018 LEAVE CREDITS" to
#     "2017 LEAVE BALANCE", with the on-screen selection left on the new
#     entries.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---------------------------------------------------------------------
# 1) Copy the formatting of the "leave usage" cells in sheet1 row 82
#    (A/B/D columns) down onto sheet2 row 18, which is about to receive
#    that same entry. This reproduces the border/number-format combo
#    (styles 39/21/40) that row 82 already had.
# ---------------------------------------------------------------------
$ws1.Range("A82:D82").Copy()
$ws2.Range("A18:D18").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Sheet2 ("2017 LEAVE BALANCE") / Table1: move the SL(2-0-0) and
#    SL(4-0-0) leave entries into rows 18-19, and add the two new
#    VL(8-0-0) entries in rows 20-21.
# ---------------------------------------------------------------------
$ws2.Range("A18").Value = 45108
$ws2.Range("B18").Value = "SL(2-0-0)"
$ws2.Range("H18").Value = 2
$ws2.Range("K18").Value = "7/3,4/2023"

$ws2.Range("B19").Value = "SL(4-0-0)"
$ws2.Range("H19").Value = 4
$ws2.Range("K19").Value = "6/26-30/2023"

$ws2.Range("A20").Value = 45170
$ws2.Range("B20").Value = "VL(8-0-0)"
$ws2.Range("D20").Value = 8
$ws2.Range("K20").Value = "9/7,13,14,19,20,26,27,28/2023"

$ws2.Range("A21").Value = 45200
$ws2.Range("B21").Value = "VL(8-0-0)"
$ws2.Range("D21").Value = 8
$ws2.Range("K21").Value = "10/3,11,12,19,20,24,25,26/2023"

# ---------------------------------------------------------------------
# 3) Sheet1 ("2018 LEAVE CREDITS") / Table13: clear out the leave-usage
#    annotations that used to sit alongside the monthly earn rows, and
#    extend the monthly earn calculator (1.25 days/month) through
#    Aug-Nov 2023 (rows 83-86).
# ---------------------------------------------------------------------
$ws1.Range("B82").Value = "SP(1-0-0)"
$ws1.Range("C82").Value = 1.25
$ws1.Range("H82").ClearContents()
$ws1.Range("K82").Value = 45132
$ws1.Range("K82").NumberFormat = "m/d/yyyy"

$ws1.Range("A83").Value = 45139
$ws1.Range("B83").ClearContents()
$ws1.Range("C83").Value = 1.25
$ws1.Range("H83").ClearContents()
$ws1.Range("K83").ClearContents()

$ws1.Range("A84").Value = 45170
$ws1.Range("B84").ClearContents()
$ws1.Range("C84").Value = 1.25
$ws1.Range("K84").Value = ""
$ws1.Range("K84").NumberFormat = "General"

$ws1.Range("A85").Value = 45200

$ws1.Range("A86").Value = 45231

# ---------------------------------------------------------------------
# 4) Delete the two now-blank trailing rows (125 & 126) from Table13 so
#    the table shrinks from A8:K127 down to A8:K125 (the former row 127
#    - the table's bottom border row - shifts up to become row 125).
# ---------------------------------------------------------------------
$ws1.Range("A125:K126").EntireRow.Delete()

# ---------------------------------------------------------------------
# 5) Leave the UI focused on the new entries on the "2017 LEAVE
#    BALANCE" sheet, matching where the author ended up after typing
#    them in.
# ---------------------------------------------------------------------
$ws1.Range("F91").Select()
$ws2.Activate()
$ws2.Range("N24").Select()
